$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 533

# Update row 3 values
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 182

# Remove the old rows 4 and 5 (last two control-point cases)
$ws.Range("A4:B5").Delete()
